$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "calendly link" header + URL for Nick's row and all Yoda rows ---
$ws.Range("E1").Value = "calendly link"
$ws.Range("E2").Value = "https://calendly.com/nick-griffiths-22/strategy-meeting-clone"
$ws.Range("E3").Value = "https://calendly.com/nick-griffiths-22/strategy-meeting-clone"

# give E1 the same (bold header) format as the rest of row 1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Columns.Item(5).ColumnWidth = 53.8

# --- Duplicate the "Yoda" row (row 3) into rows 4-6, preserving formatting & hyperlink ---
$ws.Range("A3:E3").Copy($ws.Range("A4"))
$ws.Range("A3:E3").Copy($ws.Range("A5"))
$ws.Range("A3:E3").Copy($ws.Range("A6"))

# give each duplicate a distinct name
$ws.Range("A4").Value = "Yoda 2"
$ws.Range("A5").Value = "Yoda 3"
$ws.Range("A6").Value = "Yoda 4"

# add the mailto hyperlinks for the new rows (same target as row 3)
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:yoda@email.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:yoda@email.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:yoda@email.com")

# re-apply D3's exact cell format to the new hyperlink cells so no new style is introduced
$ws.Range("D3").Copy()
$ws.Range("D4:D6").PasteSpecial(-4122)

# update selection to match the author's final cursor position
$ws.Range("A6").Select()
